$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings referenced by the new rows.
$sNadi    = "نادي الشباب"
$sSiyanaM = "صيانة المركز"
$sSaound  = "الساوند قطية"
$sDiyafaM = "ضيافة معا"
$sLafita  = "لافتة معا الجديدة"
$sDiyafah = "ضيافه"
$sSiyana  = "صيانة"

# New expense rows (rows 30..38), mirroring the style of the existing data rows.
$rows = @(
    @(29, 45666, $sNadi,    $sNadi,    1, 1000000, 20025846445),
    @(30, 45666, $sSiyanaM, $sSiyana,  1, 10000,   20024197699),
    @(31, 45666, $sSiyanaM, $sSiyana,  1, 30000,   20024244396),
    @(32, 45666, $sSaound,  $sSaound,  1, 120000,  20024243209),
    @(33, 45666, $sDiyafaM, $sDiyafah, 1, 25000,   20024244994),
    @(34, 45666, $sLafita,  $sLafita,  1, 330000,  20024503051),
    @(35, 45666, $sSiyanaM, $sSiyana,  1, 1000000, 20023537607),
    @(36, 45666, $sSiyanaM, $sSiyana,  1, 300000,  20023647624),
    @(37, 45666, $sSiyanaM, $sSiyana,  1, 500000,  20023765391)
)

$r = 30
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item(29, 1).Style

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).Style = $ws.Cells.Item(29, 2).Style

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).Style = $ws.Cells.Item(29, 1).Style

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item(29, 1).Style

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).Style = $ws.Cells.Item(29, 1).Style

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).Style = $ws.Cells.Item(29, 1).Style

    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).Style = $ws.Cells.Item(29, 1).Style

    $r = $r + 1
}

# Clear the old total row (row 30 originally held F30=SUM(F2:F29)); the new
# total moves to row 39.
$ws.Range("F39").Formula = "=SUM(F2:F38)"
$ws.Range("F39").Style = "Comma"
$ws.Range("F39").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"

$ws.Range("A13").Select()
$ws.Range("F30").Select()
